# Reorder the "Recorded By" names/emails in column G.
# The author re-ordered the comma-separated list of recorders in several
# rows (System/admin/dnasr281/backdoor variants) without changing the set
# of values present in each cell - just their order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    switch ($val) {
        "System, dnasr281@gmail.com" { $cell.Value = "dnasr281@gmail.com, System" }
        "admin@admin.com, System" { $cell.Value = "System, admin@admin.com" }
        "admin@admin.com, dnasr281@gmail.com" { $cell.Value = "dnasr281@gmail.com, admin@admin.com" }
        "backup@backdoor.com, system, System" { $cell.Value = "backup@backdoor.com, System, system" }
        default { }
    }
}
